$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the first data row (row 2), shifting all
# existing price history rows down by one.
$ws.Rows("2:2").Insert()

# The freshly inserted row inherits formatting from the row above
# (the bold header row); strip that back to the plain, unstyled look
# used by every other data row.
$ws.Rows("2:2").ClearFormats()

# Write the new day's data. The leading apostrophe keeps the date as
# literal text (matching every other date cell in column A) instead of
# letting Excel coerce it into a date serial number.
$ws.Range("A2").Value = "'2026-01-20"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Remove any residual cell style picked up from the quote-prefixed
# text entry so the row matches the unstyled cells below it.
$ws.Range("A2:D2").ClearFormats()
